$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.659.11'
$ws.Range("E2").Value = '  +2.34%  '
$ws.Range("D3").Value = '3.194.57'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.40'
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.88'
$ws.Range("E6").Value = '  +5.25%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.194.66'
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.551'
$ws.Range("E9").Value = '  +3.41%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.93'
$ws.Range("E11").Value = '  -4.17%  '
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.98'
$ws.Range("E14").Value = '  +3.67%  '
$ws.Range("D15").Value = '3.720.33'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '66.673.36'
$ws.Range("E16").Value = '  +2.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.41'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").Value = '3.195.56'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.111'
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '515.27'
$ws.Range("E20").Value = '  +1.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.39'
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.737'
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.16'
$ws.Range("E23").Value = '  +3.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.96'
$ws.Range("E24").Value = '  -0.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.01'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.27'
$ws.Range("E27").Value = '  +1.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.01'
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.41'
$ws.Range("E29").Value = '  +9.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.09'
$ws.Range("E30").Value = '  +8.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.04'
$ws.Range("E31").Value = '  +8.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.16'
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.56'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '517.81'
$ws.Range("E36").Value = '  +9.09%  '
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0897'
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0422'
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.89'
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.125'
$ws.Range("E41").Value = '  +6.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.89'
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("D43").Value = '0.0₃0686'
$ws.Range("E43").Value = '  +12.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.304'
$ws.Range("E44").Value = '  +6.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.47'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("D46").Value = '2.868.42'
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.59'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("E48").Value = '  +5.76%  '
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.65'
$ws.Range("E51").Value = '  +9.38%  '
